# Generate Report for Handoff
# Localization status has moved from "Handed back: in sync with en-US"
# to "Ready for handoff" and the corresponding timestamps advance.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language status + latest handoff-xliff-generate date ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-05 09:22:20"

# --- zh-cn sheet: status + latest handoff datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-05 09:22:15"

# --- de-de sheet: status + latest handoff datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-05 09:22:20"

# --- Column widths shrink now that the status text is shorter ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.26
$wsOverview.Columns.Item(6).ColumnWidth = 16.26
$wsZhCn.Columns.Item(3).ColumnWidth = 16.26
$wsDeDe.Columns.Item(3).ColumnWidth = 16.26
